# Refresh the stock-screener figures on Sheet1 (price/change/yield/stochastics)
# for rows 2-30, matching the newly uploaded source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D20 and D28 previously used a bespoke "0%" (no-decimals) percentage style
# while every other cell in column D uses "0.00%". Re-apply the common format
# first so, after the value refresh below, the whole column is formatted
# consistently and no cell is left pointing at the now-unused style.
$ws.Range("D20").NumberFormat = "0.00%"
$ws.Range("D28").NumberFormat = "0.00%"

$updates = @(
    @{ Row = 2; C = 22950; D = -0.0497; I = 4.36; J = 77; K = 77 },
    @{ Row = 3; C = 100800; D = -0.0147; I = 6.45; J = 64; K = 64 },
    @{ Row = 4; C = 418500; D = -0.0499; I = 4.54; J = 70; K = 70 },
    @{ Row = 5; C = 29600; D = -0.0793; I = 6.76; J = 42; K = 42 },
    @{ Row = 6; C = 30750; D = -0.0191; I = 3.9; J = 79; K = 79 },
    @{ Row = 7; C = 23650; D = -0.0744; I = 5.07; J = 62; K = 62 },
    @{ Row = 8; C = 10290; D = -0.0347; I = 5; J = 81; K = 81 },
    @{ Row = 9; C = 83500; D = -0.0402; I = 3.59; J = 74; K = 74 },
    @{ Row = 10; C = 210000; D = -0.0141; I = 5.71; J = 40; K = 40 },
    @{ Row = 11; C = 122500; D = -0.043; I = 5.55; J = 77; K = 77 },
    @{ Row = 12; C = 18660; D = -0.0762; I = 5.09; J = 68; K = 68 },
    @{ Row = 13; C = 66900; D = -0.0551; I = 5.23; J = 73; K = 73 },
    @{ Row = 14; C = 55500; D = -0.0107; I = 6.38; J = 71; K = 71 },
    @{ Row = 15; C = 83800; D = -0.0176; I = 6.56; J = 89; K = 89 },
    @{ Row = 16; C = 18650; D = -0.0396; I = 5.71; J = 74; K = 74 },
    @{ Row = 17; C = 48950; D = -0.0171; I = 5.72; J = 66; K = 66 },
    @{ Row = 18; C = 19570; D = -0.0156; I = 6.29; J = 31; K = 31 },
    @{ Row = 19; C = 54000; D = -0.0217; I = 3.7; J = 86; K = 86 },
    @{ Row = 20; C = 14470; D = -0.0109; I = 4.49; J = 74; K = 74 },
    @{ Row = 21; C = 129100; D = -0.0122; I = 4.18; J = 78; K = 78 },
    @{ Row = 22; C = 41300; D = -0.0384; I = 3.52; J = 44; K = 44 },
    @{ Row = 23; C = 65100; D = -0.0426; I = 3.32; J = 82; K = 82 },
    @{ Row = 24; C = 46800; D = -0.0507; I = 5.77; J = 61; K = 61 },
    @{ Row = 25; C = 82300; D = -0.0363; I = 4.37; J = 79; K = 79 },
    @{ Row = 26; C = 106000; D = -0.0442; I = 2.99; J = 77; K = 77 },
    @{ Row = 27; C = 13810; D = -0.0535; I = 4.71; J = 80; K = 80 },
    @{ Row = 28; C = 13450; D = -0.0303; I = 3.72; J = 79; K = 79 },
    @{ Row = 29; C = 21750; D = -0.0584; I = 4.57; J = 80; K = 80 },
    @{ Row = 30; C = 23600; D = -0.0445; I = 5.08; J = 82; K = 82 }
)

foreach ($u in $updates) {
    $ws.Range("C" + $u.Row).Value = $u.C
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("I" + $u.Row).Value = $u.I
    $ws.Range("J" + $u.Row).Value = $u.J
    $ws.Range("K" + $u.Row).Value = $u.K
}
